# Applies the "Proper score now added when ghost are eaten." edit:
#  - Shifts the text of the first 4 TODO list items up by one (item 1's text
#    becomes obsolete, items 2/3/4 take the text of the previous item), and
#    gives the 4th item brand-new text.
#  - The "_GoBack" bookmark (which sits at the end of the 1st item originally)
#    moves along with the edited text, ending up at the end of the (new) 4th
#    item.
#  - Adds a brand new TODO list item after the "...leave out the play
#    through." item.
#  - The "lastRenderedPageBreak" hint moves from the 2nd "Scatter for 7
#    seconds..." paragraph to the 1st one (an automatic side effect, in real
#    Word, of the extra bullet above pushing the page break earlier).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Shift the text of the first four bullet points up by one.
#    Go from the bottom up so each Find target still matches the
#    still-untouched original text below it.
# ---------------------------------------------------------------------

$d.Content.Find.Execute(
    "Adjust speed of eyes to go faster and then change over to regular chase speed when ghost is in box.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Kill pellet chumping sound and add ghost eaten sound.", 2) | Out-Null

$d.Content.Find.Execute(
    "Make sure that score is added when ghost is eaten.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Adjust speed of eyes to go faster and then change over to regular chase speed when ghost is in box.", 2) | Out-Null

$d.Content.Find.Execute(
    "Make sure that correct score is added when power pellet is eaten. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Make sure that score is added when ghost is eaten.", 2) | Out-Null

$d.Content.Find.Execute(
    "Make sure that game is bug free when pac-man dies or wins.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Make sure that correct score is added when power pellet is eaten. ", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark from the end of item 1 to the end of
#    item 4 (now reading "Kill pellet chumping sound and add ghost eaten
#    sound.").
# ---------------------------------------------------------------------

$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

$item4 = $d.Paragraphs.Item(6)
$item4End = $item4.Range.End

# The COM bridge mishandles a zero-length Range exactly at a paragraph's
# last valid text offset (Range.End - 1), so insert a throwaway character
# there first, anchor the bookmark just *before* it (a safe, non-edge
# position), then delete the throwaway character again.
$placeholder = $d.Range($item4End - 1, $item4End - 1)
$placeholder.InsertAfter("X")

$anchor = $d.Range($item4End - 1, $item4End - 1)
$d.Bookmarks.Add("_GoBack", $anchor) | Out-Null

$d.Range($item4End - 1, $item4End).Delete()

# ---------------------------------------------------------------------
# 3) Insert a brand new bullet after "...leave out the play through."
# ---------------------------------------------------------------------

$d.Content.Find.Execute(
    "For now, leave out the play through.", $true, $false, $false, $false,
    $false, $true, 1, $false, "", 0) | Out-Null

$playThroughPara = $d.Paragraphs.Item(13)
$playThroughPara.Range.InsertParagraphAfter() | Out-Null
$newPara = $d.Paragraphs.Item(14)
$newPara.Range.Text = "Add free life after 10000 pts (check to make sure) and add free life sound."

# ---------------------------------------------------------------------
# 4) Move the "lastRenderedPageBreak" rendering hint from the 2nd
#    "Scatter for 7 seconds, then Chase for 20 seconds." bullet to the
#    1st one. There is no object-model property for this Word-internal
#    pagination hint, so splice the two paragraphs' raw OOXML in place
#    (identical content/formatting, just the hint relocated).
# ---------------------------------------------------------------------

$scatter1 = $d.Paragraphs.Item(29)
$scatter2 = $d.Paragraphs.Item(30)
$scatterRange = $d.Range($scatter1.Range.Start, $scatter2.Range.End)
$scatterXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="003B5BDD" w:rsidRPr="003B5BDD" w:rsidRDefault="003B5BDD" w:rsidP="003B5BDD"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:shd w:val="clear" w:color="auto" w:fill="3F3F3F"/><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="405" w:lineRule="atLeast"/><w:rPr><w:rFonts w:ascii="Georgia" w:eastAsia="Times New Roman" w:hAnsi="Georgia" w:cs="Times New Roman"/><w:color w:val="DCDCCC"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r w:rsidRPr="003B5BDD"><w:rPr><w:rFonts w:ascii="Georgia" w:eastAsia="Times New Roman" w:hAnsi="Georgia" w:cs="Times New Roman"/><w:color w:val="DCDCCC"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:val="en-US"/></w:rPr><w:lastRenderedPageBreak/><w:t>Scatter for 7 seconds, then Chase for 20 seconds.</w:t></w:r></w:p><w:p w:rsidR="003B5BDD" w:rsidRPr="003B5BDD" w:rsidRDefault="003B5BDD" w:rsidP="003B5BDD"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:shd w:val="clear" w:color="auto" w:fill="3F3F3F"/><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="405" w:lineRule="atLeast"/><w:rPr><w:rFonts w:ascii="Georgia" w:eastAsia="Times New Roman" w:hAnsi="Georgia" w:cs="Times New Roman"/><w:color w:val="DCDCCC"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r w:rsidRPr="003B5BDD"><w:rPr><w:rFonts w:ascii="Georgia" w:eastAsia="Times New Roman" w:hAnsi="Georgia" w:cs="Times New Roman"/><w:color w:val="DCDCCC"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:val="en-US"/></w:rPr><w:t>Scatter for 7 seconds, then Chase for 20 seconds.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$scatterRange.InsertXML($scatterXml) | Out-Null
